$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$full = $tr.Text

foreach ($label in @("DS", "GB", "SG")) {
    $idx = $full.IndexOf("`t" + $label)
    $startPos = $idx + 1 + 1   # 1-based index of the char right after the tab
    $sub = $tr.Characters($startPos, $label.Length)
    $sub.Font.Highlight.RGB = 65535
}
